# Append the new trade row (row 6) to the GILD bag-trade log.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = 42649.644687499997
$ws.Range("B6").Value = $false
$ws.Range("C6").Value = 10000.67
$ws.Range("D6").Value = 10014.69
$ws.Range("E6").Value = 77.349997999999999
$ws.Range("F6").Value = 77.569999999999993
$ws.Range("G6").Value = $true
$ws.Range("H6").Value = 0.28000000000000003
$ws.Range("I6").Value = $false

# Match the date-formatted style already used by column A / G (style index 1)
# by copying formatting from an existing data row instead of authoring a new
# number format (which would bloat the style table).
$ws.Range("A3").Copy()
$ws.Range("A6").PasteSpecial(-4122)
$ws.Range("G3").Copy()
$ws.Range("G6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New BuyPrice values widened column E's best-fit width.
$ws.Columns.Item(5).ColumnWidth = 9
